# "3tp haste calculation adjusted"
#
# Adds a status column (done/todo) next to the existing "Functionality"
# table, adds four new functionality rows (itemID/link grabbing + OsF 3tp
# support + proper rotation support), moves the old "TODO late"-style rows
# further down the sheet to make room, gives the long wrapped label its own
# style, and tweaks the page setup / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: push the trailing block of rows (old 18-21) down to 45-48 by
# inserting 27 blank rows starting at row 8.
$ws.Range("8:34").Insert()

# New functionality rows 8-9 (column B) --------------------------------------
$ws.Range("B8").Value = "Grab itemID/link from mouseover"
$ws.Range("B9").Value = "Grab IitemID/link from mouseover comparison(shift)"

# The second one is long, so it gets a wrapped-text style + taller row.
$ws.Range("B9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 30

# Status column C for the existing functionality rows ------------------------
$ws.Range("C3").Value = "done"
$ws.Range("C4").Value = "done"
$ws.Range("C5").Value = "done"
$ws.Range("C6").Value = "done"
$ws.Range("C7").Value = "todo"

# More new functionality rows 10-11 (column B) --------------------------------
$ws.Range("B10").Value = "OsF 3tp support"
$ws.Range("B11").Value = "proper rotation support"

# Status column C for the new functionality rows ------------------------------
$ws.Range("C8").Value = "todo"
$ws.Range("C9").Value = "todo"
$ws.Range("C10").Value = "todo"
$ws.Range("C11").Value = "todo"

# Page setup ------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection ---------------------------------------------------------------
$ws.Range("C12").Select()
